$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The Price (D) / Volume (E) columns hold plain text in the workbook (the site
# renders raw scraped strings, not numbers). Most new values are already
# non-numeric-looking (two "." thousands separators, or padded "  +x.xx%  "
# strings) so a normal assignment keeps them as text. A handful of new Price
# values now parse as plain decimals (e.g. "570.79"), so a leading apostrophe
# forces Excel to keep them as text instead of silently converting to numbers.
$ws.Range("D2").Value = "64.042.53"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "3.400.91"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'570.79"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'162.52"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.400.86"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").Value = "'0.549"
$ws.Range("E9").Value = "  -4.55%  "
$ws.Range("D10").Value = "'7.30"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("E12").Value = "  -4.09%  "
$ws.Range("D13").Value = "3.991.33"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "'26.82"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "64.054.74"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "3.360.59"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").Value = "'6.10"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "'13.44"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "'372.71"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "'7.77"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'70.19"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("E25").Value = "  -3.15%  "
$ws.Range("E26").Value = "  -4.48%  "
$ws.Range("D27").Value = "'9.48"
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'6.07"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "'1.39"
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("D32").Value = "'2.00"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'22.73"
$ws.Range("E34").Value = "  -2.19%  "
$ws.Range("D35").Value = "'7.01"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("E36").Value = "  -6.02%  "
$ws.Range("D37").Value = "'159.71"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "'0.854"
$ws.Range("E38").Value = "  +7.94%  "
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").Value = "'0.0725"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("D41").Value = "'25.74"
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").Value = "'42.75"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").Value = "'6.45"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "2.724.53"
$ws.Range("E44").Value = "  -5.53%  "
$ws.Range("D45").Value = "'25.78"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").Value = "  -3.93%  "
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").Value = "'2.41"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").Value = "'330.14"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("E51").Value = "  -1.84%  "
